# Move the "Comment" column (currently column J) to the end of the
# table (after "Trp", i.e. to column AJ), shifting columns K:AJ one
# position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 36   # AJ
$lastRow = 5

# Capture the full J1:J5 "Comment" column values before we shift anything.
$commentValues = @()
for ($r = 1; $r -le $lastRow; $r++) {
    $commentValues += , ($ws.Cells.Item($r, 10).Value2)
}

# Shift columns K (11) .. AJ (36) left by one, into J (10) .. AI (35).
for ($c = 11; $c -le $lastCol; $c++) {
    for ($r = 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, $c - 1).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Place the saved Comment column values into the last column (AJ).
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $lastCol).Value = $commentValues[$r - 1]
}
